$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("B18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("B19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("B20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("B21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
